$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data: restore the pre-merge demo credentials row ---
$ws.Range("A2").Value = "https://rsepesyrcnew.rmx.rentmanager.qa/"
$ws.Range("B2").Value = "admin"
$ws.Range("C2").Value = "123!@#aA"
$ws.Range("D2").Value = "Ryan Sepesy"

# --- Hyperlink on A2 should point at the restored URL ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://rsepesyrcnew.rmx.rentmanager.qa/")

# --- Cell A2 loses the auto-applied "Hyperlink" style (back to Normal) ---
$ws.Range("A2").Style = "Normal"

# --- Selection moves to B2 ---
$ws.Range("B2").Select()
